# Apply the "nuevos experimentos no convexos" update: refresh the numeric
# experiment data (expressions / evaluations / restriction values) across the
# non-convex MitsosBarton2006Ex323 generator workbook's sheets.
#
# Every affected cell in this workbook stores its content as TEXT (shared
# string), even when that text looks like a plain number (e.g. "0.68").
# Excel's normal Range.Value assignment auto-detects numeric-looking strings
# and stores them as real numbers, which would produce a different cell type
# than the original file uses. To force a text cell (without leaving any
# permanent "Text" number-format styling behind), we:
#   1. set NumberFormat = "@" (Text) so the assignment is not re-interpreted,
#   2. assign the literal string to .Value,
#   3. reset the cell's Style back to "Normal" so no stray formatting
#      remains on the cell itself.
function Set-TextValue {
    param(
        $Worksheet,
        [string]$Address,
        [string]$Text
    )
    $rng = $Worksheet.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Text
    $rng.Style = "Normal"
}

$wb = $excel.ActiveWorkbook

# NOTE: worksheet lookup by name is case-insensitive (as in real Excel), and
# this workbook has two sheets whose names differ only by case
# ("Vector_bf" vs "Vector_BF"). Looking either of those up by name would
# resolve to whichever one appears first in tab order, silently editing the
# wrong sheet. Index into $wb.Worksheets by its (1-based) tab position
# instead, which is unambiguous:
#   1 Funciones_Objetivo
#   2 Restricciones_del_lider
#   3 Restricciones_del_follower
#   4 Punto_modificado
#   5 Vector_bf
#   6 Vector_BF
#   7 Vector_Alpha

# --- Sheet: Restricciones_del_lider ---------------------------------------
$ws = $wb.Worksheets.Item(2)
Set-TextValue $ws "A2" "1.9399999999999995 - x"
Set-TextValue $ws "B2" "-2.9399999999999995"
Set-TextValue $ws "D2" "0.32"
Set-TextValue $ws "A3" "-1.9399999999999997 + x"
Set-TextValue $ws "B3" "0.9399999999999997"
Set-TextValue $ws "D3" "0.02"
Set-TextValue $ws "A4" "35.63239999999999 + x - y - 9(x^2)"
Set-TextValue $ws "B4" "-34.63239999999999"
Set-TextValue $ws "D4" "0.44"

# --- Sheet: Restricciones_del_follower -------------------------------------
$ws = $wb.Worksheets.Item(3)
Set-TextValue $ws "A2" "-19.71359999999999 + (-0.5 + x)*(y^2)"
Set-TextValue $ws "B2" "19.71359999999999"
Set-TextValue $ws "D2" "0.36"
Set-TextValue $ws "E2" "2.6"
Set-TextValue $ws "F2" "6.7"
Set-TextValue $ws "A3" "-3.6999999999999993 + y"
Set-TextValue $ws "B3" "2.6999999999999993"
Set-TextValue $ws "D3" "0.43"
Set-TextValue $ws "E3" "7.0"
Set-TextValue $ws "F3" "0.7000000000000001"
Set-TextValue $ws "A4" "-5.699999999999999 - y"
Set-TextValue $ws "B4" "-4.699999999999999"
Set-TextValue $ws "D4" "0.43"
Set-TextValue $ws "E4" "0.7000000000000001"
Set-TextValue $ws "F4" "4.6000000000000005"

# --- Sheet: Punto_modificado -------------------------------------------
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws "A2" "1.9399999999999997"
Set-TextValue $ws "B2" "3.6999999999999993"

# --- Sheet: Vector_bf -------------------------------------------------
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws "A2" "-4.836159999999998"

# --- Sheet: Vector_BF -------------------------------------------------
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws "A2" "-24.249199999999988"
Set-TextValue $ws "A3" "-33.56559999999999"

Write-Host "Applied nuevos experimentos no convexos updates"
